$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset all cell contents (keeps formatting/styles, e.g. the bold header row)
# so the shared-string table is rebuilt from scratch in write order.
$ws.Cells.ClearContents()

# --- Header row (row 1): strings 0-19, identical to the original file ---
$ws.Range("A1").Value = "Sending cluster"
$ws.Range("B1").Value = "Ligand symbol"
$ws.Range("C1").Value = "Receptor symbol"
$ws.Range("D1").Value = "Target cluster"
$ws.Range("E1").Value = "Ligand-expressing cells"
$ws.Range("F1").Value = "Ligand detection rate"
$ws.Range("G1").Value = "Ligand average expression value"
$ws.Range("H1").Value = "Ligand total expression value"
$ws.Range("I1").Value = "Ligand derived specificity of average expression value"
$ws.Range("J1").Value = "Ligand derived specificity of total expression value"
$ws.Range("K1").Value = "Receptor-expressing cells"
$ws.Range("L1").Value = "Receptor detection rate"
$ws.Range("M1").Value = "Receptor average expression value"
$ws.Range("N1").Value = "Receptor total expression value"
$ws.Range("O1").Value = "Receptor derived specificity of average expression value"
$ws.Range("P1").Value = "Receptor derived specificity of total expression value"
$ws.Range("Q1").Value = "Edge average expression weight"
$ws.Range("R1").Value = "Edge total expression weight"
$ws.Range("S1").Value = "Edge average expression derived specificity"
$ws.Range("T1").Value = "Edge total expression derived specificity"

# --- Seed the cluster-name strings in the exact order the target workbook
#     expects them to land in the shared-string table: ECs, FAPs, Itgb2, Thy1, MuSCs.
#     (Shared-string allocation order follows write order, not sheet position,
#     so we write them once to a scratch cell before the real data goes in.)
$ws.Range("Z1").Value = "ECs"
$ws.Range("Z1").Value = "FAPs"
$ws.Range("Z1").Value = "Itgb2"
$ws.Range("Z1").Value = "Thy1"
$ws.Range("Z1").Value = "MuSCs"
$ws.Range("Z1").ClearContents()

# --- Data rows 2-7 ---
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Itgb2"
$ws.Range("C2").Value = "Thy1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.05619066666666667
$ws.Range("H2").Value = 0.168572
$ws.Range("I2").Value = 0.3931387525216601
$ws.Range("J2").Value = 0.39313875252166
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.943736666666666
$ws.Range("N2").Value = 5.83121
$ws.Range("O2").Value = 0.02216753253531823
$ws.Range("P2").Value = 0.02216753253531823
$ws.Range("Q2").Value = 0.1092198591244444
$ws.Range("R2").Value = 0.98297873212
$ws.Range("S2").Value = 0.00871491608741832
$ws.Range("T2").Value = 0.008714916087418319

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Itgb2"
$ws.Range("C3").Value = "Thy1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.05619066666666667
$ws.Range("H3").Value = 0.168572
$ws.Range("I3").Value = 0.3931387525216601
$ws.Range("J3").Value = 0.39313875252166
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 63.45677799999999
$ws.Range("N3").Value = 190.370334
$ws.Range("O3").Value = 0.7236989531682786
$ws.Range("P3").Value = 0.7236989531682786
$ws.Range("Q3").Value = 3.565678660338666
$ws.Range("R3").Value = 32.09110794304799
$ws.Range("S3").Value = 0.2845141036498083
$ws.Range("T3").Value = 0.2845141036498083

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Itgb2"
$ws.Range("C4").Value = "Thy1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.05619066666666667
$ws.Range("H4").Value = 0.168572
$ws.Range("I4").Value = 0.3931387525216601
$ws.Range("J4").Value = 0.39313875252166
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.28342866666667
$ws.Range("N4").Value = 66.850286
$ws.Range("O4").Value = 0.2541335142964031
$ws.Range("P4").Value = 0.2541335142964031
$ws.Range("Q4").Value = 1.252120712399111
$ws.Range("R4").Value = 11.269086411592
$ws.Range("S4").Value = 0.09990973278443337
$ws.Range("T4").Value = 0.09990973278443335

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Itgb2"
$ws.Range("C5").Value = "Thy1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.08673766666666667
$ws.Range("H5").Value = 0.260213
$ws.Range("I5").Value = 0.60686124747834
$ws.Range("J5").Value = 0.60686124747834
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.943736666666666
$ws.Range("N5").Value = 5.83121
$ws.Range("O5").Value = 0.02216753253531823
$ws.Range("P5").Value = 0.02216753253531823
$ws.Range("Q5").Value = 0.1685951830811111
$ws.Range("R5").Value = 1.51735664773
$ws.Range("S5").Value = 0.01345261644789991
$ws.Range("T5").Value = 0.01345261644789991

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Itgb2"
$ws.Range("C6").Value = "Thy1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.08673766666666667
$ws.Range("H6").Value = 0.260213
$ws.Range("I6").Value = 0.60686124747834
$ws.Range("J6").Value = 0.60686124747834
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 63.45677799999999
$ws.Range("N6").Value = 190.370334
$ws.Range("O6").Value = 0.7236989531682786
$ws.Range("P6").Value = 0.7236989531682786
$ws.Range("Q6").Value = 5.504092857904666
$ws.Range("R6").Value = 49.536835721142
$ws.Range("S6").Value = 0.4391848495184703
$ws.Range("T6").Value = 0.4391848495184703

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Itgb2"
$ws.Range("C7").Value = "Thy1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.08673766666666667
$ws.Range("H7").Value = 0.260213
$ws.Range("I7").Value = 0.60686124747834
$ws.Range("J7").Value = 0.60686124747834
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.28342866666667
$ws.Range("N7").Value = 66.850286
$ws.Range("O7").Value = 0.2541335142964031
$ws.Range("P7").Value = 0.2541335142964031
$ws.Range("Q7").Value = 1.932812607879778
$ws.Range("R7").Value = 17.395313470918
$ws.Range("S7").Value = 0.1542237815119697
$ws.Range("T7").Value = 0.1542237815119697
